$wb = $excel.ActiveWorkbook

# Week 16 brought a new player (J.Johnson) into the Ravens' Yards Data log.
# His column is inserted right after T.Huntley (before T.Williams) on both
# the "Rushing" and "Receiving" tracking sheets, with the same "n" sentinel
# seeded in the Week-2 (simulated) row as every other player already has.
foreach ($ws in $wb.Worksheets) {
    $ws.Columns("D").Insert()
    $ws.Range("D1").Value = "J.Johnson"
    $ws.Range("D2").Value = "n"
}
